# Updating price policy code
# The "Electric_boiler" technology column is removed from the
# Installation / Capacity results sheets, and the re-run model's new
# capacity numbers (and the Storage_capacity figure) are written in.

$wb = $excel.ActiveWorkbook

# 1. Remove the "Electric_boiler" column (column B) from the
#    Installation and Capacity sheets. This shifts every later column
#    (Gas_CHP, Gas_boiler, Grid, Heat_pump, Solar_PV, Solar_thermal)
#    one place to the left and drops the now-unused shared string.
$installation = $wb.Worksheets.Item("Installation")
$installation.Columns.Item(2).Delete()

$capacity = $wb.Worksheets.Item("Capacity")
$capacity.Columns.Item(2).Delete()

# 2. Write the updated capacity results (new model run) into the
#    Capacity sheet.
$capacity.Range("B2").Value = 11.560693641618496
$capacity.Range("C2").Value = 0
$capacity.Range("D2").Value = 0
$capacity.Range("E2").Value = 0
$capacity.Range("F2").Value = 1140.0784121337804
$capacity.Range("G2").Value = 0

$capacity.Range("B3").Value = 20
$capacity.Range("C3").Value = 563.40070657927242
$capacity.Range("D3").Value = 0
$capacity.Range("E3").Value = 0
$capacity.Range("F3").Value = 0
$capacity.Range("G3").Value = 0

# 3. Update the Storage_capacity result.
$storage = $wb.Worksheets.Item("Storage_capacity")
$storage.Range("B2").Value = 897.60282631708935
